$d = $word.ActiveDocument
$tbl = $d.Tables.Item(1)
$values = @(
  "28+59=",
  "14+79=",
  "49+33=",
  "49+29=",
  "76-57=",
  "51-39=",
  "9+84=",
  "86-48=",
  "39+23=",
  "6+89=",
  "77-38=",
  "39+23=",
  "69+7=",
  "91-45=",
  "90-46=",
  "40-24=",
  "39+42=",
  "93-68=",
  "35+38=",
  "3+29=",
  "46+46=",
  "31-7=",
  "65+18=",
  "24-15=",
  "37+8=",
  "80-17=",
  "18+54=",
  "92-69=",
  "53-7=",
  "9+37=",
  "59+19=",
  "50-6=",
  "49+32=",
  "26-17=",
  "72-43=",
  "4+27=",
  "71-17=",
  "44-15=",
  "25-9=",
  "19+13=",
  "52-48=",
  "33-24=",
  "94-45=",
  "82-3=",
  "82-8=",
  "70-27=",
  "79+15=",
  "19+77=",
  "47+29=",
  "19+36=",
  "29+57=",
  "49+8=",
  "19+32=",
  "9+6=",
  "51-7=",
  "33-17=",
  "36-19=",
  "22+49=",
  "90-15=",
  "20-6=",
  "85-46=",
  "8+75=",
  "6+79=",
  "9+73=",
  "9+78=",
  "64+18=",
  "26+5=",
  "42+39=",
  "25+7=",
  "11-3=",
  "89+7=",
  "17+24=",
  "72-44=",
  "19+64=",
  "36+55=",
  "58+9=",
  "37+18=",
  "35+37=",
  "6+27=",
  "90-67=",
  "54+29=",
  "48+9=",
  "47+34=",
  "85-68=",
  "24+68=",
  "84-55=",
  "15+18=",
  "28+4=",
  "38+47=",
  "57+7=",
  "32-26=",
  "31-22=",
  "30-7=",
  "4+19=",
  "18+23=",
  "40-33=",
  "39+29=",
  "26+49=",
  "83-28=",
  "57+7="
)

$rows = $tbl.Rows.Count
$cols = $tbl.Columns.Count

if ($rows * $cols -ne $values.Length) {
    Write-Output ("WARNING: table has " + ($rows * $cols) + " cells but " + $values.Length + " values were supplied")
}

$idx = 0
for ($r = 1; $r -le $rows; $r++) {
    for ($c = 1; $c -le $cols; $c++) {
        if ($idx -lt $values.Length) {
            $cell = $tbl.Cell($r, $c)
            $cell.Range.Text = $values[$idx]
        }
        $idx = $idx + 1
    }
}
Write-Output ("Updated " + $idx + " cells")
